# feat: add 2022-Q1 data
#
# The existing last sheet "总计" (grand-total summary) is repurposed to
# become the new "2022-Q1" per-fund detail sheet (it keeps sheetId=6 /
# sheet6.xml), and a brand-new "总计" sheet is appended after it
# (sheetId=7 / sheet7.xml) holding the refreshed summary table (the old
# summary rows shifted down by one, with a new 2022-Q1 row on top).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0) Structural changes first (sheet references become stale/re-indexed
#    once Worksheets.Add() runs, so do all add/rename/move work up front
#    and only re-fetch fresh sheet handles afterwards for the data fill).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

$q1Ref = $wb.Worksheets.Item("2022-Q1")
$newTotal.Move($null, $q1Ref)

# Fresh handles, fetched after all structural moves are done.
$q1 = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Fill the "2022-Q1" fund-detail sheet
# ---------------------------------------------------------------------

# Header row (B1:D1 already carry the bold/border header style; reuse it
# for the newly introduced E1:H1 header cells via a format-only copy).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund rows. Columns B-G hold text (fund codes have leading zeros, and the
# decimal figures are stored as text in this workbook), so force a text
# number format before assigning so the engine doesn't re-parse them as
# numbers. Column A (row index) reuses the existing bold/border style,
# column H (rank) is a plain number.
$q1.Range("B2:G8").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "011834"
$q1.Range("C2").Value = "大成投资严选六个月持有期混合型证券投资基金A"
$q1.Range("D2").Value = "3.88"
$q1.Range("E2").Value = "84.63"
$q1.Range("F2").Value = "6.27"
$q1.Range("G2").Value = "0.2433"
$q1.Range("H2").Value = 6

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "013463"
$q1.Range("C3").Value = "大成致远优势一年持有期混合A"
$q1.Range("D3").Value = "4.01"
$q1.Range("E3").Value = "60.15"
$q1.Range("F3").Value = "4.05"
$q1.Range("G3").Value = "0.1624"
$q1.Range("H3").Value = 6

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "007107"
$q1.Range("C4").Value = "太平 MSCI 香港价值增强指数A"
$q1.Range("D4").Value = "1.05"
$q1.Range("E4").Value = "93.78"
$q1.Range("F4").Value = "2.32"
$q1.Range("G4").Value = "0.0244"
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "011835"
$q1.Range("C5").Value = "大成投资严选六个月持有期混合型证券投资基金C"
$q1.Range("D5").Value = "0.30"
$q1.Range("E5").Value = "84.63"
$q1.Range("F5").Value = "6.27"
$q1.Range("G5").Value = "0.0188"
$q1.Range("H5").Value = 6

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "162416"
$q1.Range("C6").Value = "华宝港股通恒生香港35指数(LOF)"
$q1.Range("D6").Value = "0.21"
$q1.Range("E6").Value = "94.50"
$q1.Range("F6").Value = "4.36"
$q1.Range("G6").Value = "0.0092"
$q1.Range("H6").Value = 8

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "013464"
$q1.Range("C7").Value = "大成致远优势一年持有期混合C"
$q1.Range("D7").Value = "0.17"
$q1.Range("E7").Value = "60.15"
$q1.Range("F7").Value = "4.05"
$q1.Range("G7").Value = "0.0069"
$q1.Range("H7").Value = 6

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "007108"
$q1.Range("C8").Value = "太平 MSCI 香港价值增强指数C"
$q1.Range("D8").Value = "0.00"
$q1.Range("E8").Value = "93.78"
$q1.Range("F8").Value = "2.32"
# G8 is the one market-value cell that is a genuine number (0), not text.
$q1.Range("G8").NumberFormat = "General"
$q1.Range("G8").Value = 0
$q1.Range("H8").Value = 9

# Style column A (row index cells) the same as the existing A2 cell.
$q1.Range("A2").Copy()
$q1.Range("A3:A8").PasteSpecial(-4122)
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3
$q1.Range("A6").Value = 4
$q1.Range("A7").Value = 5
$q1.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2) Fill the new "总计" sheet with the refreshed summary table
# ---------------------------------------------------------------------
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.46

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 2.71

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.04

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.02

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 1
$total.Range("D7").Value = 0.01

# Apply the header/index style (copied from the 2022-Q1 sheet, which still
# carries the original style definitions) to the new sheet's header row
# and index column.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
